$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells in column D whose new value would otherwise be auto-parsed
# by Excel as a pure number need to be forced back to Text so they keep
# matching the source data's string type (e.g. "213.63", "24.25", ...).
$textCells = @("D5","D8","D10","D11","D14","D15","D16","D18","D20","D23","D25","D26","D27","D34","D36","D42","D43","D46","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.402.16"
$ws.Range("E2").Value = "  -2.33%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.651.05"
$ws.Range("E3").Value = "  -2.37%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.18%  "

# Row 5 - BNB
$ws.Range("D5").Value = "213.63"
$ws.Range("E5").Value = "  -1.59%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.18%  "

# Row 8 - Solana
$ws.Range("D8").Value = "24.25"
$ws.Range("E8").Value = "  +0.31%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.81%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0616"
$ws.Range("E10").Value = "  -1.86%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0878"
$ws.Range("E11").Value = "  -0.90%  "

# Row 12 - now WrappedEther (was WrappedliquidstakedEther2.0)
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.765.79"
$ws.Range("E12").Value = "  +4.42%  "

# Row 13 - now WrappedliquidstakedEther2.0 (was WrappedEther)
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.884.75"
$ws.Range("E13").Value = "  -2.38%  "

# Row 14 - now Polygon (was Polkadot)
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.574"
$ws.Range("E14").Value = "  +3.04%  "

# Row 15 - now Polkadot (was Polygon)
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "4.09"
$ws.Range("E15").Value = "  -2.28%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "65.99"
$ws.Range("E16").Value = "  -1.38%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.406.39"
$ws.Range("E17").Value = "  -2.20%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "234.45"
$ws.Range("E18").Value = "  -6.32%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0726"
$ws.Range("E19").Value = "  -2.28%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  -3.00%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.24%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.21%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "9.31"
$ws.Range("E23").Value = "  -2.57%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.31%  "

# Row 25 - Monero
$ws.Range("D25").Value = "147.00"
$ws.Range("E25").Value = "  -0.50%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "7.21"
$ws.Range("E26").Value = "  -1.83%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "16.09"
$ws.Range("E27").Value = "  -2.48%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  -0.07%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -2.13%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -1.40%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.89%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -2.11%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.459.64"
$ws.Range("E33").Value = "  +0.69%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "3.11"
$ws.Range("E34").Value = "  -2.51%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -4.16%  "

# Row 36 - HuobiToken
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  -0.87%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  -3.97%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -3.22%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -1.47%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  -0.30%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.18%  "

# Row 42 - now Aave (was FraxShare)
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "65.49"
$ws.Range("E42").Value = "  -5.78%  "

# Row 43 - now FraxShare (was Aave)
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  -2.33%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  -1.01%  "

# Row 45 - now RocketPoolETH (was TrustWalletToken)
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.793.64"
$ws.Range("E45").Value = "  -2.37%  "

# Row 46 - now TrustWalletToken (was RocketPoolETH)
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "0.785"
$ws.Range("E46").Value = "  -1.55%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -0.19%  "

# Row 48 - Quant
$ws.Range("E48").Value = "  -1.22%  "

# Row 49 - BabyDogeCoin
$ws.Range("D49").Value = "0.0₆0106"
$ws.Range("E49").Value = "  -4.55%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -1.71%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "7.80"
$ws.Range("E51").Value = "  -1.66%  "
